$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of (cellRef, newValue) pairs extracted from the target diff.
$changes = @(
    ,@('D2', '65.422.74')
    ,@('E2', '  -1.31%  ')
    ,@('D3', '3.431.45')
    ,@('E3', '  -4.31%  ')
    ,@('E4', '  -0.03%  ')
    ,@('D5', '594.11')
    ,@('E5', '  -1.98%  ')
    ,@('D6', '134.48')
    ,@('E6', '  -9.06%  ')
    ,@('D7', '3.430.21')
    ,@('E7', '  -4.30%  ')
    ,@('E8', '  -0.26%  ')
    ,@('E9', '  +0.29%  ')
    ,@('D10', '7.49')
    ,@('E10', '  -4.18%  ')
    ,@('D11', '0.122')
    ,@('E12', '  -8.78%  ')
    ,@('D13', '4.008.01')
    ,@('E13', '  -4.47%  ')
    ,@('D14', '0.0000180')
    ,@('E14', '  -12.59%  ')
    ,@('D15', '26.40')
    ,@('E15', '  -10.57%  ')
    ,@('D16', '65.347.43')
    ,@('E16', '  -1.55%  ')
    ,@('D17', '3.417.60')
    ,@('E17', '  -4.60%  ')
    ,@('E18', '  -2.25%  ')
    ,@('D19', '9.87')
    ,@('E19', '  -10.52%  ')
    ,@('D20', '5.73')
    ,@('E20', '  -9.43%  ')
    ,@('D21', '13.66')
    ,@('D22', '391.05')
    ,@('E22', '  -7.54%  ')
    ,@('D23', '73.08')
    ,@('E23', '  -6.83%  ')
    ,@('E24', '  -11.01%  ')
    ,@('E25', '  -0.03%  ')
    ,@('D26', '3.571.36')
    ,@('E26', '  -4.35%  ')
    ,@('E27', '  -12.76%  ')
    ,@('E28', '  -0.10%  ')
    ,@('B29', 'RenderToken')
    ,@('C29', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr')
    ,@('D29', '7.21')
    ,@('E29', '  -12.30%  ')
    ,@('B30', 'PancakeSwap')
    ,@('C30', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake')
    ,@('D30', '2.25')
    ,@('E30', '  -9.55%  ')
    ,@('E31', '  -12.83%  ')
    ,@('D32', '3.434.66')
    ,@('E32', '  -4.13%  ')
    ,@('E33', '  -0.04%  ')
    ,@('D34', '0.144')
    ,@('E34', '  -7.98%  ')
    ,@('D35', '22.65')
    ,@('E35', '  -9.45%  ')
    ,@('D36', '171.86')
    ,@('E36', '  -1.75%  ')
    ,@('E37', '  -14.43%  ')
    ,@('D38', '6.83')
    ,@('E38', '  -11.75%  ')
    ,@('E39', '  -9.18%  ')
    ,@('D40', '4.80')
    ,@('E40', '  -13.68%  ')
    ,@('D41', '0.0770')
    ,@('E41', '  -9.61%  ')
    ,@('D42', '0.811')
    ,@('E42', '  -7.79%  ')
    ,@('D43', '43.41')
    ,@('E43', '  -5.29%  ')
    ,@('E44', '  +0.02%  ')
    ,@('D45', '4.39')
    ,@('E45', '  -15.11%  ')
    ,@('E46', '  -12.82%  ')
    ,@('D47', '1.09')
    ,@('E47', '  -3.01%  ')
    ,@('D48', '22.25')
    ,@('E48', '  -5.79%  ')
    ,@('E49', '  -8.91%  ')
    ,@('E50', '  -15.94%  ')
    ,@('D51', '2.184.09')
    ,@('E51', '  -8.49%  ')
)

foreach ($pair in $changes) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $cell = $ws.Range($cellRef)
    # Force text interpretation so Excel does not coerce numeric-looking
    # strings (e.g. "0.122", "0.0000180") into actual numbers, which would
    # lose trailing zeros / formatting present in the source data.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    # Restore the default/general style so no stray explicit style index
    # is left attached to the cell (keeps cell formatting unchanged).
    $cell.Style = "Normal"
}
